# "Updated symbol list" run: refresh the Price column with new quotes,
# and roll the coin listing (rows 14-26) up by one slot - the coin that
# used to be on row 14 reappears at the bottom (row 26) with a fresh price.
#
# D-column prices are stored as TEXT in this sheet (e.g. "23.09"), not
# numbers, so each numeric-looking value is entered with a leading
# apostrophe to force text entry, then the style is reset back to
# "Normal" so no formatting residue (e.g. a "Text" number format) is
# left behind on the cell - matching the original look & feel exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "'23.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'6.384"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.06285"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Value = "'6.743"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'1.389"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.8376"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1629"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.08395"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03494"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03137"
$ws.Range("D13").Style = "Normal"
$ws.Range("B14").Value = "MCDex"
$ws.Range("C14").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D14").Value = "'3.971"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "13MCDexMCB"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09308"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001703"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04852"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Value = "'0.006282"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D19").Value = "'0.005486"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").Value = "'0.001089"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D21").Value = "'0.0001498"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "20NitroExNTX"
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").Value = "'3.737"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "21LEOLEO"
$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D23").Value = "'2.356"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "22BTSETokenBTSE"
$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D24").Value = "'0.01389"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "23OneONE"
$ws.Range("B25").Value = "BitpandaEcosystemToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D25").Value = "'0.3406"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "24BitpandaEcosystemTokenBEST"
$ws.Range("B26").Value = "ProBitToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D26").Value = "'0.1264"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "25ProBitTokenPROB"
$ws.Range("D27").Value = "'0.0002665"
$ws.Range("D27").Style = "Normal"
$ws.Range("D40").Value = "'0.04690"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.006889"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.1177"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.003451"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Value = "'0.00006255"
$ws.Range("D45").Style = "Normal"
$ws.Range("D48").Value = "'0.09913"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002097"
$ws.Range("D49").Style = "Normal"
